$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "-"
$ws.Range("D10").Value = "-"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "-"
$ws.Range("E12").Value = "MEC-3A-Usin. CNC"
$ws.Range("E14").Value = "MEC-3A-Usin. CNC"
$ws.Range("E15").Value = "MEC-3A-Usin. CNC"
$ws.Range("E16").Value = "MEC-3A-Usin. CNC"
